$wb = $excel.ActiveWorkbook

# --- Rubric score update on "Student Score" sheet ---
# Item 8 ("View not responding to Post parameters") actual score goes from 0 -> 2.
# The Subtotal (C26 = SUM(C24:C25)) and Total (C28 = SUM(C16,C21,C26)) cells are
# formulas, so they recalculate automatically from this single edit.
$scoreSheet = $wb.Worksheets.Item("Student Score")
$scoreSheet.Range("C25").Value = 2

# --- View-state changes ---
# Previously "Lab3Rubric_CS295N" was the selected/active tab with E4 selected;
# now "Student Score" is the active tab (with E20 selected) and
# "Lab3Rubric_CS295N" is left with A4 selected.
$rubricSheet = $wb.Worksheets.Item("Lab3Rubric_CS295N")
$rubricSheet.Range("A4").Select()

$scoreSheet.Activate()
$scoreSheet.Range("E20").Select()
